# Fix: Make Signature Status Check robust using Server Time (updated_at) instead of Client Time

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "cliente" value in K2 from "Brinks" to "Carozzi"
$ws.Range("K2").Value = "Carozzi"

# Update ticket_id value in B2 from 12 to 13
$ws.Range("B2").Value = 13

# Update the active selection / active cell to I6
$ws.Range("I6").Select()
